$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" column (strikeouts) values replacing the old "Strike#" counts.
$kValues = @{
    2  = 3
    3  = 6
    4  = 2
    5  = 5
    6  = 4
    7  = 1
    8  = 2
    9  = 8
    10 = 6
    11 = 7
    12 = 2
    13 = 4
    14 = 5
    15 = 2
    16 = 4
    17 = 2
    18 = 4
    19 = 4
    20 = 5
    21 = 4
    22 = 4
    23 = 7
    24 = 1
    25 = 7
    26 = 6
    27 = 2
    28 = 3
    29 = 8
    30 = 2
    31 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
